# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Tuna" (Vega Modelo de Temuco) as row 88,
# shifting the existing rows 88-93 down to 89-94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 88 (pushes 88..93 -> 89..94)
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row with the new data record
$ws.Range("A88").Value2 = 10
$ws.Range("B88").Value2 = "Vega Modelo de Temuco"
$ws.Range("C88").Value2 = "La Araucanía"
$ws.Range("D88").Value2 = 45041
$ws.Range("E88").Value2 = 9
$ws.Range("F88").Value2 = "Fruta"
$ws.Range("G88").Value2 = 100107
$ws.Range("H88").Value2 = "Otros"
$ws.Range("I88").Value2 = 100107011
$ws.Range("J88").Value2 = "Tuna"
$ws.Range("K88").Value2 = "Sin especificar"
$ws.Range("L88").Value2 = "Primera"
$ws.Range("M88").Value2 = 25
$ws.Range("N88").Value2 = 22000
$ws.Range("O88").Value2 = 22000
$ws.Range("P88").Value2 = 22000
$ws.Range("Q88").Value2 = "$/caja 16 kilos"
$ws.Range("R88").Value2 = "Provincia de Los Andes"
$ws.Range("S88").Value2 = 1375
$ws.Range("T88").Value2 = 16
